# Adding new test case for Font check
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("hub")

# Insert two new columns before column H (shifts H:M -> J:O)
$ws.Range("H1:I1").EntireColumn.Insert()

# New header cells for the inserted columns
$ws.Range("H1").Value = "width"
$ws.Range("I1").Value = "hight"

# New "70px" values for Search (row 5) and Settings (row 6)
$ws.Range("H5").Value = "70px"
$ws.Range("I5").Value = "70px"
$ws.Range("H6").Value = "70px"
$ws.Range("I6").Value = "70px"

# The old font-size cells for Search/Settings rows moved to column L; clear them
$ws.Range("L5").ClearContents()
$ws.Range("L6").ClearContents()

# Match the new column widths of the inserted columns to column G (12.140625)
$ws.Range("H1").ColumnWidth = $ws.Range("G1").ColumnWidth
$ws.Range("I1").ColumnWidth = $ws.Range("G1").ColumnWidth

# Update the recorded selection to match the authored view state
$ws.Range("I25").Select()
